$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 384
$ws.Range("I41").Value = 362.625
$ws.Range("J41").Value = 426.75
$ws.Range("K41").Value = 362.625
$ws.Range("L41").Value = 426.75
$ws.Range("M41").Value = 77.375

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2615.2307
$ws.Range("I70").Value = 4000
$ws.Range("J70").Value = 1999.7778
$ws.Range("K70").Value = 12000
$ws.Range("L70").Value = 5999.3334
$ws.Range("M70").Value = -11730

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 2615.2307
$ws.Range("I73").Value = 4000
$ws.Range("J73").Value = 1999.7778
$ws.Range("K73").Value = 12000
$ws.Range("L73").Value = 5999.3334
$ws.Range("M73").Value = -11064

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 311.8
$ws.Range("I103").Value = 549.5
$ws.Range("J103").Value = 275.23077
$ws.Range("K103").Value = 1648.5
$ws.Range("L103").Value = 825.69231
$ws.Range("M103").Value = -1062.5
$ws.Range("N103").Value = -1997.69231

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 10360.714
$ws.Range("I129").Value = 11905.4
$ws.Range("J129").Value = 6499
$ws.Range("K129").Value = 35716.2
$ws.Range("L129").Value = 19497
$ws.Range("M129").Value = -30716.2
$ws.Range("N129").Value = -29497

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2719.7273
$ws.Range("I137").Value = 1989.75
$ws.Range("J137").Value = 4666.3335
$ws.Range("K137").Value = 5969.25
$ws.Range("L137").Value = 13999.0005
$ws.Range("M137").Value = -3419.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H140").Value = 94402
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 94402
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 94402
$ws.Range("N140").Value = -104762

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3284.6296
$ws.Range("I122").Value = 2129.35
$ws.Range("J122").Value = 6585.4287
$ws.Range("K122").Value = 6388.049999999999
$ws.Range("L122").Value = 19756.2861
$ws.Range("M122").Value = -3938.049999999999
$ws.Range("N122").Value = -24656.2861

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4518.269
$ws.Range("I132").Value = 2562.0527
$ws.Range("J132").Value = 9828
$ws.Range("K132").Value = 7686.158100000001
$ws.Range("L132").Value = 29484
$ws.Range("M132").Value = -5156.158100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1404.6
$ws.Range("I86").Value = 1219.381
$ws.Range("J86").Value = 1836.7778
$ws.Range("K86").Value = 1219.381
$ws.Range("L86").Value = 1836.7778
$ws.Range("M86").Value = -96.38100000000009
$ws.Range("N86").Value = -4082.7778

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1404.6
$ws.Range("I89").Value = 1219.381
$ws.Range("J89").Value = 1836.7778
$ws.Range("K89").Value = 6096.905000000001
$ws.Range("L89").Value = 9183.889000000001
$ws.Range("M89").Value = -480.9050000000007
$ws.Range("N89").Value = -20415.889

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 785.94116
$ws.Range("I94").Value = 657.5333000000001
$ws.Range("J94").Value = 1749
$ws.Range("K94").Value = 657.5333000000001
$ws.Range("L94").Value = 1749
$ws.Range("M94").Value = -206.5333000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3166.4211
$ws.Range("I99").Value = 3041.4375
$ws.Range("J99").Value = 3833
$ws.Range("K99").Value = 3041.4375
$ws.Range("L99").Value = 3833
$ws.Range("M99").Value = -1543.4375
$ws.Range("N99").Value = -6829

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1942.2
$ws.Range("I107").Value = 1693.2354
$ws.Range("J107").Value = 2471.25
$ws.Range("K107").Value = 1693.2354
$ws.Range("L107").Value = 2471.25
$ws.Range("M107").Value = 226.7646

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6330.364
$ws.Range("I31").Value = 2645.7273
$ws.Range("J31").Value = 13699.637
$ws.Range("K31").Value = 2645.7273
$ws.Range("L31").Value = 13699.637
$ws.Range("M31").Value = -2350.7273

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 6330.364
$ws.Range("I34").Value = 2645.7273
$ws.Range("J34").Value = 13699.637
$ws.Range("K34").Value = 2645.7273
$ws.Range("L34").Value = 13699.637
$ws.Range("M34").Value = -2443.7273

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H45").Value = 50000
$ws.Range("I45").Value = 50000
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 50000
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -49407
$ws.Range("N45").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3400.9092
$ws.Range("I122").Value = 3059.6875
$ws.Range("J122").Value = 4310.8335
$ws.Range("K122").Value = 9179.0625
$ws.Range("L122").Value = 12932.5005
$ws.Range("M122").Value = -6729.0625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 5588.1113
$ws.Range("I62").Value = 3261
$ws.Range("J62").Value = 10242.333
$ws.Range("K62").Value = 9783
$ws.Range("L62").Value = 30726.999
$ws.Range("M62").Value = -9097
$ws.Range("N62").Value = -32098.999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 27250
$ws.Range("I63").Value = 27250
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 81750
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -81001
$ws.Range("N63").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 8589.6
$ws.Range("I64").Value = 7999.6
$ws.Range("J64").Value = 9179.6
$ws.Range("K64").Value = 23998.8
$ws.Range("L64").Value = 27538.8
$ws.Range("M64").Value = -23728.8
$ws.Range("N64").Value = -28078.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H65").Value = 5588.1113
$ws.Range("I65").Value = 3261
$ws.Range("J65").Value = 10242.333
$ws.Range("K65").Value = 29349
$ws.Range("L65").Value = 92180.997
$ws.Range("M65").Value = -25917
$ws.Range("N65").Value = -99044.997

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H66").Value = 27250
$ws.Range("I66").Value = 27250
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 245250
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -241506
$ws.Range("N66").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H67").Value = 8589.6
$ws.Range("I67").Value = 7999.6
$ws.Range("J67").Value = 9179.6
$ws.Range("K67").Value = 23998.8
$ws.Range("L67").Value = 27538.8
$ws.Range("M67").Value = -23062.8
$ws.Range("N67").Value = -29410.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2996.5
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 2996.5
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 8989.5
$ws.Range("N68").Value = -10611.5
$ws.Range("M68").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 2996.5
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 2996.5
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 26968.5
$ws.Range("N71").Value = -35080.5
$ws.Range("M71").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 20466.334
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 20466.334
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 61399.00199999999
$ws.Range("N74").Value = -63521.00199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 3406
$ws.Range("I75").Value = 1671.4
$ws.Range("J75").Value = 4645
$ws.Range("K75").Value = 5014.200000000001
$ws.Range("L75").Value = 13935
$ws.Range("M75").Value = -4016.200000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H77").Value = 20466.334
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 20466.334
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 184197.006
$ws.Range("N77").Value = -194805.006

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H78").Value = 3406
$ws.Range("I78").Value = 1671.4
$ws.Range("J78").Value = 4645
$ws.Range("K78").Value = 15042.6
$ws.Range("L78").Value = 41805
$ws.Range("M78").Value = -10050.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 11117.577
$ws.Range("I80").Value = 9146.471
$ws.Range("J80").Value = 14840.777
$ws.Range("K80").Value = 9146.471
$ws.Range("L80").Value = 14840.777
$ws.Range("M80").Value = -8148.471
$ws.Range("N80").Value = -16836.777

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 11117.577
$ws.Range("I83").Value = 9146.471
$ws.Range("J83").Value = 14840.777
$ws.Range("K83").Value = 45732.355
$ws.Range("L83").Value = 74203.88499999999
$ws.Range("M83").Value = -40740.355
$ws.Range("N83").Value = -84187.88499999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 46963.547
$ws.Range("I7").Value = 61044.24
$ws.Range("J7").Value = 17394.1
$ws.Range("K7").Value = 61044.24
$ws.Range("L7").Value = 17394.1
$ws.Range("M7").Value = -60932.24

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 2966.6667
$ws.Range("I32").Value = 2966.6667
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 2966.6667
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -2649.6667
$ws.Range("N32").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3597.5264
$ws.Range("I46").Value = 943
$ws.Range("J46").Value = 5528.091
$ws.Range("K46").Value = 943
$ws.Range("L46").Value = 5528.091
$ws.Range("M46").Value = -755
$ws.Range("N46").Value = -5904.091

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1102.4
$ws.Range("I82").Value = 1126
$ws.Range("J82").Value = 1067
$ws.Range("K82").Value = 1126
$ws.Range("L82").Value = 1067
$ws.Range("M82").Value = -765

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1102.4
$ws.Range("I85").Value = 1126
$ws.Range("J85").Value = 1067
$ws.Range("K85").Value = 1126
$ws.Range("L85").Value = 1067
$ws.Range("M85").Value = 122

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5659.5
$ws.Range("I122").Value = 4638.5
$ws.Range("J122").Value = 9403.166999999999
$ws.Range("K122").Value = 13915.5
$ws.Range("L122").Value = 28209.501
$ws.Range("M122").Value = -11465.5
$ws.Range("N122").Value = -33109.501

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 46963.547
$ws.Range("I126").Value = 61044.24
$ws.Range("J126").Value = 17394.1
$ws.Range("K126").Value = 183132.72
$ws.Range("L126").Value = 52182.3
$ws.Range("M126").Value = -180662.72

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 29596.555
$ws.Range("I96").Value = 2699.5
$ws.Range("J96").Value = 37281.43
$ws.Range("K96").Value = 2699.5
$ws.Range("L96").Value = 37281.43
$ws.Range("M96").Value = -1326.5
$ws.Range("N96").Value = -40027.43

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2888.348
$ws.Range("I122").Value = 3018.4285
$ws.Range("J122").Value = 2686
$ws.Range("K122").Value = 9055.2855
$ws.Range("L122").Value = 8058
$ws.Range("M122").Value = -6605.2855
